$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# Add "(IPv6, Zero Trust, Zero-Config VPN, etc.) " before "und Backup-Strategien."
# in the "Überwacht die Konzeption..." bullet.
$null = $d.Content.Find.Execute(
    "VPN-Lösungen und Backup-Strategien.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "VPN-Lösungen (IPv6, Zero Trust, Zero-Config VPN, etc.) und Backup-Strategien.",
    2)

# --- Change 2 ---------------------------------------------------------
# Insert a new bullet "Erfolg: Implementierung von IPv6 im gesamten
# Unternehmen, Zukunftssicherung des Netzwerks." right before the
# "Erfolg: Ausrollen sicherer SD-WAN-Anbindungen..." bullet.
$rng2 = $d.Content
$null = $rng2.Find.Execute(
    "Erfolg: Ausrollen sicherer SD-WAN-Anbindungen",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$rng2.Collapse(1)  # wdCollapseStart

$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="1B6E5A"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>' +
        '<w:t xml:space="preserve">Erfolg: Implementierung von IPv6 im gesamten Unternehmen, Zukunftssicherung des Netzwerks.</w:t></w:r></w:p>'
$rng2.InsertXML($xml2)

# --- Change 3 ---------------------------------------------------------
# Insert a new bullet "Erfolg: Konzeption und Implementierung von IPv6
# im Netzwerk-Perimeter zur Sicherstellung der externen Erreichbarkeit."
# right after the "Erfolg: Aufbau einer Mail-Archivierungsplattform..."
# bullet.
$rng3 = $d.Content
$null = $rng3.Find.Execute(
    "Erfolg: Aufbau einer Mail-Archivierungsplattform, die Langzeit-Compliance sicherstellte und eDiscovery vereinfachte.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$rng3.Collapse(0)  # wdCollapseEnd

$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="1B6E5A"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>' +
        '<w:t xml:space="preserve">Erfolg: Konzeption und Implementierung von IPv6 im Netzwerk-Perimeter zur Sicherstellung der externen Erreichbarkeit.</w:t></w:r></w:p>'
$rng3.InsertXML($xml3)

$d.Save()
